# Update as of 2024-03-17
# ------------------------------------------------------------------
# This script reproduces the 2024-03-17 update to the "aportes"
# workbook:
#   - Ingreso: remove the erroneous -200 "Punto" row (old row 666),
#     which shifts everything below it up by one row, then append the
#     eight new entries recorded on 2024-03-17 (date serial 45368).
#   - Gastos: correct the date on the last "Arbitro y agua" entry and
#     add the new 2024-03-17 expense (900+140).
#   - Cuentas por cobrar: Kibelo's outstanding debt was paid, so that
#     row is removed from the receivables sheet...
#   - Histórico de tecnicas: ...and appears here instead, along with a
#     new technical-foul entry for Johan. Two pre-existing rows also
#     get their "Concepto" corrected from "Tecnica" to "Técnica".
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ===================== Ingreso (sheet1) =====================
$ws1 = $wb.Worksheets.Item("Ingreso")

# Remove the incorrect row (2024-03-03, Punto, -200) — everything
# below shifts up one row.
$ws1.Rows.Item(666).Delete()

# The autofilter / filter-database range no longer covers the new
# rows being appended below (it stays pinned at the old bottom, row
# 676) — set this *before* appending so it doesn't auto-grow.
$ws1.AutoFilterMode = $false
$ws1.Range("A1:D676").AutoFilter()
$wb.Names.Item(1).RefersTo = "=Ingreso!`$A`$1:`$D`$676"

# Append the 2024-03-17 entries.
$newRows1 = @(
    @(45368, "Johan",     100, "Técnica"),
    @(45368, "Kibelo",    100, "Técnica"),
    @(45368, "Kibelo",    900, "Aporte"),
    @(45368, "Invitados", 400, "Aporte"),
    @(45368, "Punto",     300, "Aporte"),
    @(45368, "Randy",     100, "Aporte"),
    @(45368, "Jordan",    100, "Aporte"),
    @(45368, "Johan",     300, "Aporte")
)

$r = 677
foreach ($row in $newRows1) {
    $ws1.Cells.Item($r, 1).Value = $row[0]
    $ws1.Cells.Item($r, 2).Value = $row[1]
    $ws1.Cells.Item($r, 3).Value = $row[2]
    $ws1.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

$ws1.Activate()
$ws1.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws1.Range("A685").Select()

# ===================== Gastos (sheet2) =====================
$ws2 = $wb.Worksheets.Item("Gastos")

# Fix the date on the last recorded "Arbitro y agua" expense.
$ws2.Cells.Item(84, 1).Value = 45361

# Add the new 2024-03-17 expense.
$ws2.Cells.Item(85, 1).Value = 45368
$ws2.Cells.Item(85, 2).Value = "Arbitro y agua"
$ws2.Range("C85").Formula = "=900+140"

$ws2.Activate()
$ws2.Range("A85").Select()

# ============== Cuentas por cobrar (sheet3) ==============
$ws3 = $wb.Worksheets.Item("Cuentas por cobrar")

# Kibelo's debt (2023-07-09) was paid off — remove the row; the rows
# below shift up one.
$ws3.Rows.Item(3).Delete()

$ws3.Activate()
$ws3.Range("F15").Select()

# ============== Histórico de tecnicas (sheet4) ==============
$ws4 = $wb.Worksheets.Item("Histórico de tecnicas")

# Correct the "Concepto" spelling on two existing entries.
$ws4.Range("C11").Value = "Técnica"
$ws4.Range("C12").Value = "Técnica"

# Duplicate formatting from row 12 onto the two new rows first, so
# the new date cells reuse the existing date style instead of minting
# a new one.
$ws4.Range("A12:F12").Copy()
$ws4.Range("A13:F13").PasteSpecial(-4122)
$ws4.Range("A14:F14").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New technical-foul entry for Johan (2024-03-17).
$ws4.Cells.Item(13, 1).Value = 45368
$ws4.Cells.Item(13, 2).Value = "Johan"
$ws4.Cells.Item(13, 3).Value = "Técnica"
$ws4.Cells.Item(13, 4).Value = 100
$ws4.Cells.Item(13, 5).Value = $true
$ws4.Cells.Item(13, 6).Value = "Hizo un pique con el mmg de carlos y le amagó un trompón"

# Kibelo's now-resolved receivable, moved over from "Cuentas por cobrar".
$ws4.Cells.Item(14, 1).Value = 45067
$ws4.Cells.Item(14, 2).Value = "Kibelo"
$ws4.Cells.Item(14, 3).Value = "Tecnica"
$ws4.Cells.Item(14, 4).Value = 100
$ws4.Cells.Item(14, 5).Value = $true
$ws4.Cells.Item(14, 6).Value = "Cogió un pique y picó la pelota muy duro"

$ws4.Activate()
$ws4.Range("A14").Select()

# Leave "Ingreso" as the active sheet, matching the original workbook.
$ws1.Activate()
